$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.10%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.13%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.63%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07526"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.08%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.569"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.25%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9298"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.15%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.408"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.46%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1215"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.49%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1827"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.56%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08857"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.74%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04070"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.35%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1053"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.02%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.11%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005885"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.73%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.346"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.56%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.362"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.80%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3287"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.07%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.970"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.26%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1418"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.68%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2962"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.74%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04054"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.10%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001263"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.38%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003905"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.38%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001228"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.19%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.03%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02423"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.06%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.83%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005901"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-10.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007801"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.54%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1330"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.40%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007358"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.11%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007823"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.68%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2976"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.86%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006315"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.93%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.25%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04518"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-48.34%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004193"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.20%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.25%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.25%"
